# Auto-generated Excel COM-interop script
# Applies the numeric corrections described in the commit diff
# (values refreshed via a scheduled "Sheets" runner) to each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1404.091
$ws.Range("I6").Value = 1404.091
$ws.Range("K6").Value = 4212.272999999999
$ws.Range("M6").Value = -4100.272999999999
# Row 8
$ws.Range("H8").Value = 2773.889
$ws.Range("I8").Value = 2773.889
$ws.Range("K8").Value = 8321.667000000001
$ws.Range("M8").Value = -8182.667000000001
# Row 40
$ws.Range("H40").Value = 2812.5
# Row 52
$ws.Range("H52").Value = 1320
$ws.Range("I52").Value = 600
$ws.Range("K52").Value = 1800
$ws.Range("M52").Value = -1640
# Row 86
$ws.Range("H86").Value = 1673.3
$ws.Range("I86").Value = 2198.6
$ws.Range("K86").Value = 2198.6
$ws.Range("M86").Value = -1075.6
# Row 89
$ws.Range("H89").Value = 1673.3
$ws.Range("I89").Value = 2198.6
$ws.Range("K89").Value = 10993
$ws.Range("M89").Value = -5377
# Row 98
$ws.Range("H98").Value = 2359.6956
$ws.Range("I98").Value = 2288.2778
$ws.Range("K98").Value = 2288.2778
$ws.Range("M98").Value = -790.2777999999998
# Row 112
$ws.Range("H112").Value = 2653.6177
$ws.Range("I112").Value = 1592.3334
$ws.Range("J112").Value = 2756.3225
$ws.Range("K112").Value = 4777.0002
$ws.Range("L112").Value = 8268.967500000001
$ws.Range("M112").Value = -3669.0002
$ws.Range("N112").Value = -10484.9675
# Row 122
$ws.Range("H122").Value = 2359.6956
$ws.Range("I122").Value = 2288.2778
$ws.Range("K122").Value = 6864.8334
$ws.Range("M122").Value = -4414.8334
# Row 135
$ws.Range("H135").Value = 16673349
$ws.Range("I135").Value = 20002620
$ws.Range("J135").Value = 27000
$ws.Range("K135").Value = 180023580
$ws.Range("L135").Value = 243000
$ws.Range("M135").Value = -180021045
$ws.Range("N135").Value = -248070
# Row 137
$ws.Range("H137").Value = 13895590
$ws.Range("I137").Value = 41667884
$ws.Range("K137").Value = 125003652
$ws.Range("M137").Value = -125001102

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22415.568
$ws.Range("I32").Value = 22813.96
$ws.Range("K32").Value = 22813.96
$ws.Range("M32").Value = -22526.96
# Row 102
$ws.Range("H102").Value = 4100.4707
$ws.Range("I102").Value = 3586.2856
$ws.Range("J102").Value = 6500
$ws.Range("K102").Value = 3586.2856
$ws.Range("L102").Value = 6500
$ws.Range("M102").Value = -1964.2856
$ws.Range("N102").Value = -9744
# Row 122
$ws.Range("H122").Value = 1480.6666
$ws.Range("I122").Value = 1479.5294
$ws.Range("K122").Value = 4438.5882
$ws.Range("M122").Value = -1988.5882
# Row 132
$ws.Range("H132").Value = 4530.5938
$ws.Range("I132").Value = 2856.36
$ws.Range("K132").Value = 8569.08
$ws.Range("M132").Value = -6039.08

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3647.0527
$ws.Range("I20").Value = 3302.3635
$ws.Range("K20").Value = 3302.3635
$ws.Range("M20").Value = -3055.3635
# Row 86
$ws.Range("H86").Value = 577611
$ws.Range("J86").Value = 672962.8
$ws.Range("L86").Value = 672962.8
$ws.Range("N86").Value = -675208.8
# Row 88
$ws.Range("H88").Value = 33217.75
$ws.Range("J88").Value = 33217.75
$ws.Range("L88").Value = 33217.75
$ws.Range("N88").Value = -34029.75
# Row 89
$ws.Range("H89").Value = 577611
$ws.Range("J89").Value = 672962.8
$ws.Range("L89").Value = 3364814
$ws.Range("N89").Value = -3376046
# Row 91
$ws.Range("H91").Value = 33217.75
$ws.Range("J91").Value = 33217.75
$ws.Range("L91").Value = 33217.75
$ws.Range("N91").Value = -36025.75
# Row 105
$ws.Range("H105").Value = 142865950
$ws.Range("I105").Value = 166676500
$ws.Range("K105").Value = 166676500
$ws.Range("M105").Value = -166674753
# Row 134
$ws.Range("H134").Value = 7638.364
$ws.Range("I134").Value = 4418
$ws.Range("K134").Value = 13254
$ws.Range("M134").Value = -10719

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 29
$ws.Range("H29").Value = 5608.222
$ws.Range("J29").Value = 5608.222
$ws.Range("L29").Value = 5608.222
$ws.Range("N29").Value = -6194.222
# Row 105
$ws.Range("H105").Value = 1245
$ws.Range("I105").Value = 1245
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1245
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 502
$ws.Range("N105").ClearContents()
# Row 107
$ws.Range("H107").Value = 1743.3334
$ws.Range("I107").Value = 1322.2
$ws.Range("J107").Value = 2796.1667
$ws.Range("K107").Value = 1322.2
$ws.Range("L107").Value = 2796.1667
$ws.Range("M107").Value = 597.8
$ws.Range("N107").Value = -6636.1667
# Row 122
$ws.Range("H122").Value = 92712.09
$ws.Range("J122").Value = 2072.1667
$ws.Range("L122").Value = 6216.500100000001
$ws.Range("N122").Value = -11116.5001
# Row 125
$ws.Range("H125").Value = 81267.14
$ws.Range("J125").Value = 81267.14
$ws.Range("L125").Value = 81267.14
$ws.Range("N125").Value = -86187.14
# Row 132
$ws.Range("H132").Value = 92470.21000000001
$ws.Range("I132").Value = 5053.8335
$ws.Range("J132").Value = 158032.5
$ws.Range("K132").Value = 15161.5005
$ws.Range("L132").Value = 474097.5
$ws.Range("M132").Value = -12631.5005
$ws.Range("N132").Value = -479157.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 942.75
$ws.Range("I14").Value = 942.75
$ws.Range("K14").Value = 2828.25
$ws.Range("M14").Value = -2655.25
# Row 37
$ws.Range("H37").Value = 89990
$ws.Range("J37").Value = 89990
$ws.Range("L37").Value = 269970
$ws.Range("N37").Value = -270194
# Row 69
$ws.Range("H69").Value = 1244.1428
$ws.Range("I69").Value = 1120
$ws.Range("J69").Value = 1554.5
$ws.Range("K69").Value = 3360
$ws.Range("L69").Value = 4663.5
$ws.Range("M69").Value = -2549
$ws.Range("N69").Value = -6285.5
# Row 72
$ws.Range("H72").Value = 1244.1428
$ws.Range("I72").Value = 1120
$ws.Range("J72").Value = 1554.5
$ws.Range("K72").Value = 10080
$ws.Range("L72").Value = 13990.5
$ws.Range("M72").Value = -6024
$ws.Range("N72").Value = -22102.5
# Row 92
$ws.Range("H92").Value = 1088.3214
$ws.Range("J92").Value = 905.7222
$ws.Range("L92").Value = 2717.1666
$ws.Range("N92").Value = -5213.1666
# Row 93
$ws.Range("H93").Value = 3825
$ws.Range("I93").Value = 1500
$ws.Range("K93").Value = 4500
$ws.Range("M93").Value = -2628
# Row 107
$ws.Range("H107").Value = 1955.0952
$ws.Range("I107").Value = 571.53845
$ws.Range("J107").Value = 2575.3103
$ws.Range("K107").Value = 1714.61535
$ws.Range("L107").Value = 7725.9309
$ws.Range("M107").Value = 205.38465
$ws.Range("N107").Value = -11565.9309

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Range("H53").Value = 49000.6
$ws.Range("I53").Value = 48342
$ws.Range("J53").Value = 49988.5
$ws.Range("K53").Value = 48342
$ws.Range("L53").Value = 49988.5
$ws.Range("M53").Value = -47711
$ws.Range("N53").Value = -51250.5
# Row 113
$ws.Range("H113").Value = 29498.334
$ws.Range("I113").Value = 4997
$ws.Range("K113").Value = 4997
$ws.Range("M113").Value = -2827
# Row 122
$ws.Range("H122").Value = 7425
$ws.Range("I122").Value = 8661.799999999999
$ws.Range("K122").Value = 25985.4
$ws.Range("M122").Value = -23535.4
# Row 132
$ws.Range("H132").Value = 6437.2383
$ws.Range("I132").Value = 3242.6667
$ws.Range("J132").Value = 8833.166999999999
$ws.Range("K132").Value = 9728.000100000001
$ws.Range("L132").Value = 26499.501
$ws.Range("M132").Value = -7198.000100000001
$ws.Range("N132").Value = -31559.501

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 210479
$ws.Range("J7").Value = 25497.5
$ws.Range("L7").Value = 25497.5
$ws.Range("N7").Value = -25721.5
# Row 32
$ws.Range("H32").Value = 393
$ws.Range("I32").Value = 393
$ws.Range("K32").Value = 393
$ws.Range("M32").Value = -76
# Row 50
$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 15000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 15000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -14363
$ws.Range("N50").ClearContents()
# Row 93
$ws.Range("H93").Value = 4250.5
$ws.Range("I93").Value = 4250.5
$ws.Range("K93").Value = 4250.5
$ws.Range("M93").Value = -3002.5
# Row 122
$ws.Range("H122").Value = 4303.067
$ws.Range("I122").Value = 3452.0454
$ws.Range("J122").Value = 6643.375
$ws.Range("K122").Value = 10356.1362
$ws.Range("L122").Value = 19930.125
$ws.Range("M122").Value = -7906.136200000001
$ws.Range("N122").Value = -24830.125
# Row 126
$ws.Range("H126").Value = 210479
$ws.Range("J126").Value = 25497.5
$ws.Range("L126").Value = 76492.5
$ws.Range("N126").Value = -81432.5
# Row 132
$ws.Range("H132").Value = 5847.1816
$ws.Range("I132").Value = 5017.25
$ws.Range("K132").Value = 15051.75
$ws.Range("M132").Value = -12521.75
# Row 141
$ws.Range("H141").Value = 109999
$ws.Range("J141").Value = 109999
$ws.Range("L141").Value = 109999
$ws.Range("N141").Value = -120359

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 30
$ws.Range("H30").Value = 20480.2
$ws.Range("J30").Value = 24475.25
$ws.Range("L30").Value = 24475.25
$ws.Range("N30").Value = -24689.25
# Row 38
$ws.Range("H38").Value = 24662.25
$ws.Range("I38").Value = 42499.5
$ws.Range("K38").Value = 42499.5
$ws.Range("M38").Value = -42026.5
# Row 43
$ws.Range("H43").Value = 45000
$ws.Range("J43").Value = 75000
$ws.Range("L43").Value = 75000
$ws.Range("N43").Value = -75298
# Row 133
$ws.Range("H133").Value = 40942.4
$ws.Range("J133").Value = 40942.4
$ws.Range("L133").Value = 40942.4
$ws.Range("N133").Value = -51062.4
